$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "SubArrayWithGivenSum"
$ws.Range("A10").Value = "Sub Array With Given Sum"

$ws.Range("A10").Select()
